$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.2088
$ws.Range("D9").Value = -7.496599999999994
$ws.Range("E12").Value = 17.98050000000003
$ws.Range("D18").Value = -8.939899999999993
$ws.Range("D20").Value = -7.529999999999996
$ws.Range("E26").Value = 16.1611
$ws.Range("D27").Value = -8.143499999999998
$ws.Range("E27").Value = 16.75289999999999
$ws.Range("E29").Value = 16.86540000000002
$ws.Range("E37").Value = 16.72990000000002
$ws.Range("E38").Value = 16.4146
$ws.Range("E51").Value = 17.36820000000001
$ws.Range("E55").Value = 16.4335
$ws.Range("D69").Value = -7.843599999999997
$ws.Range("E69").Value = 16.5697
$ws.Range("E70").Value = 18.00660000000002
$ws.Range("D76").Value = -7.681599999999996
$ws.Range("D82").Value = -8.237299999999991
$ws.Range("E83").Value = 16.5824
$ws.Range("E102").Value = 16.7898
